# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across
# the per-job profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 262.8
$ws.Range("I33").Value = 297.375
$ws.Range("K33").Value = 297.375
$ws.Range("M33").Value = -68.375
$ws.Range("H43").Value = 3340.9092
$ws.Range("I43").Value = 4378.2856
$ws.Range("J43").Value = 1525.5
$ws.Range("K43").Value = 4378.2856
$ws.Range("L43").Value = 1525.5
$ws.Range("M43").Value = -4309.2856
$ws.Range("N43").Value = -1663.5
$ws.Range("H98").Value = 1806.081
$ws.Range("I98").Value = 1806.25
$ws.Range("K98").Value = 1806.25
$ws.Range("M98").Value = -308.25
$ws.Range("H107").Value = 858.93335
$ws.Range("I107").Value = 877.2308
$ws.Range("K107").Value = 877.2308
$ws.Range("M107").Value = 1042.7692
$ws.Range("H122").Value = 1806.081
$ws.Range("I122").Value = 1806.25
$ws.Range("K122").Value = 5418.75
$ws.Range("M122").Value = -2968.75
$ws.Range("H135").Value = 729.1579
$ws.Range("I135").Value = 489.94446
$ws.Range("K135").Value = 4409.50014
$ws.Range("M135").Value = -1874.50014
$ws.Range("H141").Value = 4182.5
$ws.Range("I141").Value = 2453.68
$ws.Range("J141").Value = 10356.857
$ws.Range("K141").Value = 7361.039999999999
$ws.Range("L141").Value = 31070.571
$ws.Range("M141").Value = -2181.039999999999
$ws.Range("N141").Value = -41430.571
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3318.5247
$ws.Range("I32").Value = 3248.3103
$ws.Range("K32").Value = 3248.3103
$ws.Range("M32").Value = -2961.3103
$ws.Range("H45").Value = 31634.867
$ws.Range("I45").Value = 40818.727
$ws.Range("K45").Value = 40818.727
$ws.Range("M45").Value = -40441.727
$ws.Range("H74").Value = 187571.67
$ws.Range("I74").Value = 328022.75
$ws.Range("J74").Value = 3904.8462
$ws.Range("K74").Value = 328022.75
$ws.Range("L74").Value = 3904.8462
$ws.Range("M74").Value = -327148.75
$ws.Range("N74").Value = -5652.8462
$ws.Range("H77").Value = 187571.67
$ws.Range("I77").Value = 328022.75
$ws.Range("J77").Value = 3904.8462
$ws.Range("K77").Value = 1640113.75
$ws.Range("L77").Value = 19524.231
$ws.Range("M77").Value = -1635745.75
$ws.Range("N77").Value = -28260.231
$ws.Range("H122").Value = 3479.0715
$ws.Range("I122").Value = 3285.1155
$ws.Range("K122").Value = 9855.3465
$ws.Range("M122").Value = -7405.3465
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15628667
$ws.Range("I20").Value = 22731570
$ws.Range("K20").Value = 22731570
$ws.Range("M20").Value = -22731323
$ws.Range("H105").Value = 20002872
$ws.Range("I105").Value = 1252270
$ws.Range("K105").Value = 1252270
$ws.Range("M105").Value = -1250523
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1916.6666
$ws.Range("I22").Value = 1875.25
$ws.Range("K22").Value = 1875.25
$ws.Range("M22").Value = -1525.25
$ws.Range("H31").Value = 5441.6772
$ws.Range("I31").Value = 3561.6667
$ws.Range("K31").Value = 3561.6667
$ws.Range("M31").Value = -3266.6667
$ws.Range("H34").Value = 5441.6772
$ws.Range("I34").Value = 3561.6667
$ws.Range("K34").Value = 3561.6667
$ws.Range("M34").Value = -3359.6667
$ws.Range("H58").Value = 2383.6829
$ws.Range("I58").Value = 1813.7037
$ws.Range("K58").Value = 1813.7037
$ws.Range("M58").Value = -1610.7037
$ws.Range("H107").Value = 850.8570999999999
$ws.Range("I107").Value = 580.3
$ws.Range("J107").Value = 1527.25
$ws.Range("K107").Value = 580.3
$ws.Range("L107").Value = 1527.25
$ws.Range("M107").Value = 1339.7
$ws.Range("N107").Value = -5367.25
$ws.Range("H120").Value = 36957.332
$ws.Range("J120").Value = 36957.332
$ws.Range("L120").Value = 36957.332
$ws.Range("N120").Value = -44215.332
$ws.Range("H134").Value = 2191.3877
$ws.Range("J134").Value = 3207.125
$ws.Range("L134").Value = 9621.375
$ws.Range("N134").Value = -14691.375
$ws.Range("H136").Value = 2383.6829
$ws.Range("I136").Value = 1813.7037
$ws.Range("K136").Value = 5441.1111
$ws.Range("M136").Value = -2891.1111
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3474.5
$ws.Range("I80").Value = 3474
$ws.Range("J80").Value = 3475
$ws.Range("K80").Value = 10422
$ws.Range("L80").Value = 10425
$ws.Range("M80").Value = -9486
$ws.Range("N80").Value = -12297
$ws.Range("H83").Value = 3474.5
$ws.Range("I83").Value = 3474
$ws.Range("J83").Value = 3475
$ws.Range("K83").Value = 31266
$ws.Range("L83").Value = 31275
$ws.Range("M83").Value = -26586
$ws.Range("N83").Value = -40635
$ws.Range("H132").Value = 2453.15
$ws.Range("J132").Value = 3489.6
$ws.Range("L132").Value = 31406.4
$ws.Range("N132").Value = -36466.39999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9515
$ws.Range("I102").Value = 1254.4286
$ws.Range("K102").Value = 1254.4286
$ws.Range("M102").Value = 367.5714
$ws.Range("H113").Value = 2746.8572
$ws.Range("I113").Value = 2746.8572
$ws.Range("K113").Value = 2746.8572
$ws.Range("M113").Value = -576.8571999999999
$ws.Range("H126").Value = 8779.4
$ws.Range("I126").Value = 2900
$ws.Range("K126").Value = 8700
$ws.Range("M126").Value = -6230
$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3472.611
$ws.Range("J61").Value = 3329.6667
$ws.Range("L61").Value = 3329.6667
$ws.Range("N61").Value = -3733.6667
$ws.Range("H113").Value = 3472.611
$ws.Range("J113").Value = 3329.6667
$ws.Range("L113").Value = 3329.6667
$ws.Range("N113").Value = -7669.6667
$ws.Range("H140").Value = 74999.664
$ws.Range("J140").Value = 74999.664
$ws.Range("L140").Value = 74999.664
$ws.Range("N140").Value = -85359.664
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 551.86957
$ws.Range("J113").Value = 511.16666
$ws.Range("L113").Value = 1533.49998
$ws.Range("N113").Value = -5873.499980000001
$ws.Range("H122").Value = 9617571
$ws.Range("I122").Value = 2147.7368
$ws.Range("K122").Value = 6443.2104
$ws.Range("M122").Value = -3993.2104
$ws.Range("H126").Value = 1536.2727
$ws.Range("I126").Value = 1244.1111
$ws.Range("K126").Value = 3732.3333
$ws.Range("M126").Value = -1262.3333
